$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Revert "20241026 차경환 사풍 용어 추가"
# E10 previously held the "사풍" (newly-added term) explanation string;
# restore it to the original "공격력" explanation string.
$ws.Range("E10").Value = "플레이어 캐릭터의 공격력이 20 증가한다."

# Update the active sheet view/selection to match the reverted state.
$ws.Activate()
$ws.Range("E11").Select()
$excel.ActiveWindow.ScrollColumn = 2
